# "Fixed error in file naming." -- a handful of team names in the bracket
# list were mis-entered; correct them in place (values only, the cells
# keep their existing position/formatting).
#   B19: Howard      -> Wagner
#   B51: Grambling   -> Grambling St.
#   B63: Virginia    -> Colorado St.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Activate()

$ws.Range("B19").Value = "Wagner"
$ws.Range("B63").Value = "Colorado St."
$ws.Range("B51").Value = "Grambling St."

# Match the author's final view/selection state.
$ws.Application.ActiveWindow.ScrollRow = 25
$ws.Application.ActiveWindow.ScrollColumn = 1
$ws.Range("D45").Select()
